$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new test case data
$ws.Range("A2").Value = "iProc_TC_ID_44"
$ws.Range("B2").Value = "@Regression Validation of  Cloud indication -  (Status is InCompleted)"
$ws.Range("C2").Value = "passed"

# Delete row 3 (the old second test case row) entirely
$ws.Rows("3:3").Delete()
